$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "First amenment in file"
$ws.Range("A2").Select()
